{"js": "// Locate the paragraph that contains the \"{{ company }}\" merge field\n// (the \"Offerta redatta da: {{ issued_by }}, {{ company }}\" line) and\n// insert two new paragraphs right after it:\n//   1) an empty paragraph\n//   2) a paragraph containing \"{{ quotation_intro }}\"\n// Both new paragraphs inherit the same paragraph formatting as the\n// \"{{ company }}\" paragraph (Word's InsertParagraphAfter copies the\n// paragraph mark's formatting), matching the target OOXML.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst companyPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"{{ company }}\") !== -1\n);\n\nif (!companyPara) {\n  throw new Error('Could not find the paragraph containing \"{{ company }}\"');\n}\n\n// Insert the blank spacer paragraph directly after the company line.\nconst blankPara = companyPara.insertParagraph(\"\", \"After\");\n\n// Insert the quotation_intro paragraph directly after the blank one.\nblankPara.insertParagraph(\"{{ quotation_intro }}\", \"After\");\n\nawait context.sync();\n", "ps1": "# Locate the paragraph that contains the \"{{ company }}\" merge field\n# (the \"Offerta redatta da: {{ issued_by }}, {{ company }}\" line) and\n# insert two new paragraphs right after it:\n#   1) an empty paragraph\n#   2) a paragraph containing \"{{ quotation_intro }}\"\n# Both new paragraphs inherit the paragraph formatting of the\n# \"{{ company }}\" paragraph, because InsertParagraphAfter duplicates the\n# paragraph mark's formatting - matching the target OOXML.\n\n$d = $word.ActiveDocument\n\n$companyPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*{{ company }}*\") {\n        $companyPara = $p\n        break\n    }\n}\n\nif ($companyPara -eq $null) {\n    throw 'Could not find the paragraph containing \"{{ company }}\"'\n}\n\n$r = $companyPara.Range\n$r.Collapse(0)          # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n# Move into the newly created (blank) paragraph, then insert the second\n# paragraph right after it and set its text.\n$r.Collapse(0)          # wdCollapseEnd\n$r.InsertParagraphAfter()\n$r.Collapse(0)          # wdCollapseEnd\n$r.Text = \"{{ quotation_intro }}\"\n"}
